$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$win = $excel.ActiveWindow
$win.Left = 100
$win.Top = 200
$win.Width = 300
$win.Height = 400
